$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.485.96'
$ws.Range('E2').Value = '  +1.80%  '
$ws.Range('D3').Value = '3.387.23'
$ws.Range('E3').Value = '  +3.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.87%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.385.47'
$ws.Range('E8').Value = '  +3.45%  '
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.55'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.27%  '
$ws.Range('E11').Value = '  +7.43%  '
$ws.Range('E12').Value = '  +6.20%  '
$ws.Range('D13').Value = '3.954.35'
$ws.Range('E13').Value = '  +3.53%  '
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('E15').Value = '  +7.10%  '
$ws.Range('D16').Value = '3.379.02'
$ws.Range('E16').Value = '  +3.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.34'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.15%  '
$ws.Range('D18').Value = '61.472.28'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.05'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.90'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '389.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +10.27%  '
$ws.Range('E23').Value = '  +3.72%  '
$ws.Range('D24').Value = '3.519.81'
$ws.Range('E24').Value = '  +3.78%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000127'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +17.63%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '71.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.54%  '
$ws.Range('E28').Value = '  +13.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.72'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.16%  '
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.30'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.159'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.75%  '
$ws.Range('E33').Value = '  +2.90%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = '3.412.36'
$ws.Range('E35').Value = '  +3.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.51'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.35%  '
$ws.Range('E38').Value = '  +2.52%  '
$ws.Range('E39').Value = '  +4.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '161.56'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0793'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('B43').Value = 'ONDO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.23'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.68%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.72'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +11.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.771'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.71%  '
$ws.Range('E46').Value = '  +2.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '41.24'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.44'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.98'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.06'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.77%  '
$ws.Range('D51').Value = '2.366.62'
$ws.Range('E51').Value = '  +9.89%  '
